$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '69.853.48'
$ws.Range('E2').Value = '  -0.75%  '
$ws.Range('D3').Value = '3.500.03'
$ws.Range('E3').Value = '  -1.95%  '
$ws.Range('E4').Value = '  -0.29%  '
$ws.Range('D5').Value = '607.66'
$ws.Range('E5').Value = '  -0.14%  '
$ws.Range('D6').Value = '197.23'
$ws.Range('E6').Value = '  +4.95%  '
$ws.Range('D7').Value = '0.626'
$ws.Range('E7').Value = '  +0.67%  '
$ws.Range('E8').Value = '  -0.10%  '
$ws.Range('D9').Value = '0.212'
$ws.Range('E9').Value = '  -1.20%  '
$ws.Range('D10').Value = '0.658'
$ws.Range('E10').Value = '  +1.40%  '
$ws.Range('D11').Value = '54.21'
$ws.Range('E11').Value = '  +0.25%  '
$ws.Range('E12').Value = '  -0.91%  '
$ws.Range('D13').Value = '9.61'
$ws.Range('E13').Value = '  +1.94%  '
$ws.Range('D14').Value = '4.057.93'
$ws.Range('E14').Value = '  -1.91%  '
$ws.Range('D15').Value = '601.57'
$ws.Range('E15').Value = '  +4.64%  '
$ws.Range('D16').Value = '69.933.34'
$ws.Range('E16').Value = '  -0.75%  '
$ws.Range('B17').Value = 'Uniswap'
$ws.Range('C17').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D17').Value = '12.70'
$ws.Range('E17').Value = '  -0.79%  '
$ws.Range('B18').Value = 'Chainlink'
$ws.Range('C18').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D18').Value = '19.00'
$ws.Range('E18').Value = '  -0.04%  '
$ws.Range('D19').Value = '3.483.17'
$ws.Range('E19').Value = '  -2.49%  '
$ws.Range('D21').Value = '0.993'
$ws.Range('E21').Value = '  -0.45%  '
$ws.Range('D22').Value = '17.90'
$ws.Range('E22').Value = '  +2.47%  '
$ws.Range('D23').Value = '104.70'
$ws.Range('E23').Value = '  +10.85%  '
$ws.Range('D24').Value = '4.66'
$ws.Range('E24').Value = '  -2.88%  '
$ws.Range('D25').Value = '5.09'
$ws.Range('E25').Value = '  +4.37%  '
$ws.Range('D26').Value = '3.10'
$ws.Range('E26').Value = '  +5.37%  '
$ws.Range('D27').Value = '10.98'
$ws.Range('E27').Value = '  +0.16%  '
$ws.Range('D28').Value = '9.83'
$ws.Range('E28').Value = '  +4.32%  '
$ws.Range('D29').Value = '33.95'
$ws.Range('E29').Value = '  +4.87%  '
$ws.Range('D30').Value = '4.58'
$ws.Range('E30').Value = '  +23.15%  '
$ws.Range('D31').Value = '7.21'
$ws.Range('E31').Value = '  +1.63%  '
$ws.Range('D32').Value = '12.66'
$ws.Range('E32').Value = '  +3.50%  '
$ws.Range('E33').Value = '  +0.93%  '
$ws.Range('D34').Value = '64.06'
$ws.Range('E34').Value = '  -0.73%  '
$ws.Range('D35').Value = '3.683.03'
$ws.Range('E35').Value = '  -2.68%  '
$ws.Range('B36').Value = 'Bittensor'
$ws.Range('C36').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D36').Value = '522.83'
$ws.Range('E36').Value = '  +0.10%  '
$ws.Range('B37').Value = 'Dai'
$ws.Range('C37').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D37').Value = '0.999'
$ws.Range('E37').Value = '  -0.08%  '
$ws.Range('B38').Value = 'PEPE'
$ws.Range('C38').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D38').Value = '0.0₃0798'
$ws.Range('E38').Value = '  +1.79%  '
$ws.Range('B39').Value = 'Fetch.AI'
$ws.Range('C39').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D39').Value = '3.03'
$ws.Range('E39').Value = '  -5.16%  '
$ws.Range('D40').Value = '0.392'
$ws.Range('E40').Value = '  -3.46%  '
$ws.Range('D41').Value = '36.89'
$ws.Range('E41').Value = '  -3.16%  '
$ws.Range('E42').Value = '  +0.46%  '
$ws.Range('E43').Value = '  -1.25%  '
$ws.Range('D44').Value = '0.0461'
$ws.Range('E44').Value = '  +0.77%  '
$ws.Range('D45').Value = '2.86'
$ws.Range('E45').Value = '  -3.69%  '
$ws.Range('E46').Value = '  +0.16%  '
$ws.Range('E47').Value = '  -4.48%  '
$ws.Range('B48').Value = 'THORChain'
$ws.Range('C48').Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range('D48').Value = '8.77'
$ws.Range('E48').Value = '  -4.82%  '
$ws.Range('B49').Value = 'FirstDigitalUSD'
$ws.Range('C49').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D49').Value = '1.00'
$ws.Range('E49').Value = '  +0.12%  '
$ws.Range('D50').Value = '132.48'
$ws.Range('E50').Value = '  -2.94%  '
$ws.Range('E51').Value = '  -5.33%  '
